$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Edit 1: "Prerequisite(s): N/A"  ->  "Prerequisite(s): CSCI 150 or CSCI 100"
# ---------------------------------------------------------------------------
$rng1 = $d.Content
$rng1.Find.Execute("Prerequisite(s): N/A")
$text1 = $rng1.Text
$colonIdx1 = $text1.IndexOf(": N/A")
$naStart = $rng1.Start + $colonIdx1 + 2
$naEnd = $rng1.Start + $colonIdx1 + 5

# First, split the run so that "N/A" becomes its own run (same rPr as the
# colon/space run it came from). Toggling Bold on/off forces the engine to
# materialize a run boundary without leaving stray direct formatting behind.
$naRange = $d.Range($naStart, $naEnd)
$naRange.Font.Bold = $true

# Replace the (now isolated) "N/A" run's text with the new prerequisite text
# while its formatting still differs from its neighbor (keeps the engine
# from silently re-merging the two runs back together).
$naRange2 = $d.Range($naStart, $naEnd)
$naRange2.Text = "CSCI 150 or CSCI 100"

# Restore Bold back to its original (false) value on the freshly typed text.
$newLen1 = "CSCI 150 or CSCI 100".Length
$newRange1 = $d.Range($naStart, $naStart + $newLen1)
$newRange1.Font.Bold = $false

# ---------------------------------------------------------------------------
# Edit 2: "Section Number: 01"  ->  "Section Number: 00"
# ---------------------------------------------------------------------------
$rng2 = $d.Content
$rng2.Find.Execute("Section Number: 01")
$text2 = $rng2.Text
$colonIdx2 = $text2.IndexOf(": 01")
$numEnd = $rng2.Start + $colonIdx2 + 4

# Split off the trailing "1" character into its own run.
$lastCharRange = $d.Range($numEnd - 1, $numEnd)
$lastCharRange.Font.Bold = $true

$lastCharRange2 = $d.Range($numEnd - 1, $numEnd)
$lastCharRange2.Text = "0"

$newRange2 = $d.Range($numEnd - 1, $numEnd)
$newRange2.Font.Bold = $false
